# The sheet's API response columns (H,I,K,L,M) previously held shared-string
# JSON/log blobs (e.g. {"code":200,...}); the real edit replaces them with the
# plain numeric HTTP status code 200. Column J (login response) is replaced
# with the literal text of the matching "NName1"/"NName2" values already used
# in column B, reusing those shared strings instead of the old error text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 200
$ws.Range("I2").Value = 200
$ws.Range("J2").Value = "NName1"
$ws.Range("K2").Value = 200
$ws.Range("L2").Value = 200
$ws.Range("M2").Value = 200

# Row 3
$ws.Range("H3").Value = 200
$ws.Range("I3").Value = 200
$ws.Range("J3").Value = "NName2"
$ws.Range("K3").Value = 200
$ws.Range("L3").Value = 200
$ws.Range("M3").Value = 200

# The selection moves from B2 to M2:M3 (active cell M2)
$ws.Range("M2:M3").Select()
